# Update "想去人数" (want-to-go count) values on the 展览 / 演出 / 全部类型 sheets
# to reflect the latest generated numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 783
$wsExhibit.Range("F8").Value = 344
$wsExhibit.Range("F9").Value = 454
$wsExhibit.Range("F10").Value = 515
$wsExhibit.Range("F11").Value = 141
$wsExhibit.Range("F12").Value = 11772
$wsExhibit.Range("F13").Value = 5415

# 演出 (sheet2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 107

# 全部类型 (sheet4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 783
$wsAll.Range("F4").Value = 107
$wsAll.Range("F10").Value = 344
$wsAll.Range("F11").Value = 454
$wsAll.Range("F12").Value = 515
$wsAll.Range("F13").Value = 141
$wsAll.Range("F14").Value = 11772
$wsAll.Range("F16").Value = 5415
